# Timesheet_Joran_Vanhaste.xlsx update
# Fills in the timesheet table (dates, tasks, hours), formats the header,
# applies date/wrap-text styles, merges the title row, sizes columns,
# adds the totals formula and restores page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title row (A1) -------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:C1").Merge() | Out-Null

# --- Header row (row 2) ----------------------------------------------
$ws.Range("A2").Value = "Datum"
$ws.Range("B2").Value = "Taak"
$ws.Range("C2").Value = "uur"

# --- Data rows --------------------------------------------------------
# Row 3
$ws.Range("A3").Value = 42661
$ws.Range("B3").Value = "sitemap+ productiedossier opstarten"
$ws.Range("C3").Value = 2

# Row 4
$ws.Range("A4").Value = 42663
$ws.Range("B4").Value = "Moodboard + style tiles"
$ws.Range("C4").Value = 3

# Row 5
$ws.Range("A5").Value = 42664
$ws.Range("B5").Value = "Style Tiles afwerken"
$ws.Range("C5").Value = 1

# Row 6
$ws.Range("A6").Value = 42666
$ws.Range("B6").Value = "Wireframes"
$ws.Range("C6").Value = 3

# Row 7
$ws.Range("A7").Value = 42667
$ws.Range("B7").Value = "Wireframes afwerken + visual designs"
$ws.Range("C7").Value = 4

# Row 8
$ws.Range("A8").Value = 42668
$ws.Range("B8").Value = "Visuals afwerken"
$ws.Range("C8").Value = 2

# Row 9
$ws.Range("A9").Value = 42669
$ws.Range("B9").Value = "productiedossier samenstellen en afwerken, alles pushen"
$ws.Range("C9").Value = 3

# Row 11 (typed before row 10 originally)
$ws.Range("A11").Value = 42679
$ws.Range("B11").Value = "contactformulier in drupal 7 toevoegen + logo en naam veranderen"
$ws.Range("C11").Value = 1

# Row 10
$ws.Range("A10").Value = 42672
$ws.Range("B10").Value = "juiste mappenstructuur "

# Row 13 (typed before row 12 originally)
$ws.Range("A13").Value = 42682
$ws.Range("B13").Value = "Disclaimer + footer link+ Wijzigingen pushen"
$ws.Range("C13").Value = 1.5

# Row 12
$ws.Range("A12").Value = 42681
$ws.Range("B12").Value = "mappenstructuur gitlab aanpassen + repository gitlab "
$ws.Range("C12").Value = 1.5

# --- Totals -------------------------------------------------------------
$ws.Range("C14").Formula = "=SUM(C3:C13)"

# --- Number format for the date column ----------------------------------
$ws.Range("A3:A13").NumberFormat = "mm-dd-yy"

# --- Wrap text for the longer task descriptions --------------------------
$ws.Range("B3").WrapText = $true
$ws.Range("B7").WrapText = $true
$ws.Range("B9").WrapText = $true
$ws.Range("B10").WrapText = $true
$ws.Range("B11").WrapText = $true
$ws.Range("B12").WrapText = $true
$ws.Range("B13").WrapText = $true

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 7.8333333333
$ws.Columns.Item(2).ColumnWidth = 19.5

# --- Row heights (auto-computed by Excel when it wraps the text) -----------
$ws.Rows.Item(1).RowHeight = 24
$ws.Rows.Item(3).RowHeight = 48
$ws.Rows.Item(7).RowHeight = 32
$ws.Rows.Item(9).RowHeight = 48
$ws.Rows.Item(11).RowHeight = 64
$ws.Rows.Item(12).RowHeight = 48
$ws.Rows.Item(13).RowHeight = 48

# --- Selection / view state --------------------------------------------
$ws.Range("A2:C13").Select() | Out-Null

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.PaperSize = 9     # xlPaperA4
$ws.PageSetup.Orientation = 1   # xlPortrait
